# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (G) holds newly-computed s_vals for each row (2..76).
# Write the recomputed values directly into column G, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sVals = @(
    1,1,0,1,0,2,1,2,0,3,
    1,2,0,1,0,0,0,2,0,1,
    1,0,0,2,2,1,1,1,0,0,
    1,1,1,2,1,0,1,2,3,1,
    1,0,1,2,1,2,0,2,0,0,
    1,1,1,0,1,2,1,0,1,0,
    1,1,1,1,2,1,1,1,0,2,
    1,1,0,1,1
)

$firstRow = 2
for ($i = 0; $i -lt $sVals.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 7).Value = $sVals[$i]
}
